$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Fill in row 56 (12:45 - 14:00, "Remise en forme VueX + rédaction brève du rapport")
$ws.Range("C56").Value = 0.53125
$ws.Range("D56").Value = "-"
$ws.Range("E56").Value = 0.58333333333333337
$ws.Range("F56").Value = "Remise en forme VueX + rédaction brève du rapport"

# Fill in row 57 (14:00 - 15:30, "Stylisation d'editShoot  et des ArrowItems")
$ws.Range("C57").Value = 0.58333333333333337
$ws.Range("D57").Value = "-"
$ws.Range("E57").Value = 0.64583333333333337
$ws.Range("F57").Value = "Stylisation d'editShoot  et des ArrowItems"

# Update the selected cell in the sheet view
$ws.Range("C58").Select()
